$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Scene Manager progress: from "?%" to 50% (numeric, percent format)
$ws.Range("B15").Value = 0.5
$ws.Range("B15").NumberFormat = "0%"

# New shared strings must be introduced in the same order the original
# author typed them so sharedStrings.xml ordering matches (63..71):
# bang record ki luc, ?, Nam: 09/03, Nam: 09/03(Tao UI), Menu trong man choi,
# Dung game..., 1 man, ?%(Bo), Theo doi... creep ...
$ws.Range("D27").Value = "bảng record kỉ lục"
$ws.Range("B27").Value = "?"
$ws.Range("C25").Value = "Nam: 09/03"
$ws.Range("C27").Value = "Nam: 09/03(Tạo UI)"
$ws.Range("A29").Value = "Menu trong màn chơi"
$ws.Range("E29").Value = "Dừng game, chuyển ra menu chính, chỉnh nhạc, exit"
$ws.Range("D21").Value = "1 màn"
$ws.Range("B23").Value = "?%(Bỏ)"
$ws.Range("E15").Value = "Theo dõi quá trình của màn chơi, quản lý waves creep và sinh boss?"

# Remaining edits that reuse existing shared strings
$ws.Range("B25").Value = 0.5
$ws.Range("B25").NumberFormat = "0%"
$ws.Range("B28").Value = "?"
$ws.Range("B29").Value = "?"

# Old note row (old row 30) content removed entirely
$ws.Rows(30).ClearContents()

# Update selection to match the saved view state
$ws.Range("E21").Select()
